$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the more-precise timestamp value for the existing last row (20)
$ws.Range("A20").Value = 45874.79190704861

# Append the new data row (21) reported by the automatic WSL update
$ws.Range("A21").Value = 45874.83353430877
$ws.Range("A21").NumberFormat = $ws.Range("A20").NumberFormat

$ws.Range("B21").Value = 2025
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 15.78
$ws.Range("E21").Value = 86.40000000000001
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = "-"
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = "20:00:17"
